$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 7 - new subject (21.4.21)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "21.4.21"
$ws.Range("C7").Value = 24
$ws.Range("F7").Value = "left"
$ws.Range("G7").Value = "lenses"
$ws.Range("I7").Value = "F"
$ws.Range("J7").Value = "Y"
$ws.Range("L7").Value = "Y"
$ws.Range("M7").Value = "Y"
$ws.Range("N7").Value = "naomivaknine@mail.tau.ac.il"
$ws.Hyperlinks.Add($ws.Range("N7"), "mailto:naomivaknine@mail.tau.ac.il") | Out-Null

# ---------------------------------------------------------------------------
# Row 8 - new subject (21.4.21)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "21.4.21"
$ws.Range("C8").Value = 22
$ws.Range("F8").Value = "right"
$ws.Range("G8").Value = "N"
$ws.Range("I8").Value = "M"
$ws.Range("J8").Value = "Y"
$ws.Range("L8").Value = "Y"
$ws.Range("M8").Value = "Y"
$ws.Range("N8").Value = "peleg4008@gmail.com"
$ws.Hyperlinks.Add($ws.Range("N8"), "mailto:peleg4008@gmail.com") | Out-Null

# ---------------------------------------------------------------------------
# Row 9 - new subject (26.4.21)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "26.4.21"
$ws.Range("C9").Value = 23
$ws.Range("F9").Value = "right"
$ws.Range("G9").Value = "lenses"
$ws.Range("I9").Value = "F"
$ws.Range("J9").Value = "Y"
$ws.Range("K9").Value = "credit"
$ws.Range("L9").Value = "Y"
$ws.Range("M9").Value = "Y"
$ws.Range("N9").Value = "tamarsela@mail.tau.ac.il"
$ws.Hyperlinks.Add($ws.Range("N9"), "mailto:tamarsela@mail.tau.ac.il") | Out-Null

# ---------------------------------------------------------------------------
# Row 10 - new subject (26.4.21)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "26.4.21"
$ws.Range("C10").Value = 22
$ws.Range("F10").Value = "right"
$ws.Range("G10").Value = "N"
$ws.Range("I10").Value = "F"
$ws.Range("J10").Value = "Y"
$ws.Range("K10").Value = "credit"
$ws.Range("L10").Value = "Y"
$ws.Range("M10").Value = "Y"
$ws.Range("N10").Value = "rotemasher98@gmail.com"
$ws.Hyperlinks.Add($ws.Range("N10"), "mailto:rotemasher98@gmail.com") | Out-Null

# ---------------------------------------------------------------------------
# Row 11 - new subject (27.4.21)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "27.4.21"
$ws.Range("C11").Value = 25
$ws.Range("F11").Value = "left"
$ws.Range("G11").Value = "N"
$ws.Range("I11").Value = "F"
$ws.Range("J11").Value = "Y"
$ws.Range("K11").Value = "credit"
$ws.Range("L11").Value = "Y"
$ws.Range("M11").Value = "Y"
$ws.Range("N11").Value = "yuvalheimann@mail.tau.ac.il"
$ws.Hyperlinks.Add($ws.Range("N11"), "mailto:yuvalheimann@mail.tau.ac.il") | Out-Null

# ---------------------------------------------------------------------------
# Restore formatting for rows 7-11 (values above may have re-keyed some
# styles, e.g. the hyperlink cells); re-apply the canonical format from the
# already-filled-in row 5 template so every column keeps its original style.
# ---------------------------------------------------------------------------
$ws.Range("A5:N5").Copy() | Out-Null
$ws.Range("A7:N11").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Row 12 - subject slot not run yet: clear the placeholder answers but keep
# the row (B12 stays), and give N12 the same (blank) hyperlink-style cell.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "27.4.21"
$ws.Range("D12").Clear() | Out-Null
$ws.Range("E12").Clear() | Out-Null
$ws.Range("H12").Clear() | Out-Null
$ws.Range("K12").Clear() | Out-Null
$ws.Range("N5").Copy() | Out-Null
$ws.Range("N12").PasteSpecial(-4122) | Out-Null
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("A12").Style = $ws.Range("A11").Style

# ---------------------------------------------------------------------------
# Row 25 - new entry row (27.4.21 / subject 1012), same notes as row 24
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "27.4.21"
$ws.Range("B25").Value = 1012
$ws.Range("T25").Value = "Khen heller"
$ws.Range("U25").Value = "Full run on myself with full arm reach to check timing on diff screen (Asus, refrate 100) "
$ws.Rows.Item(25).RowHeight = 15.75

# ---------------------------------------------------------------------------
# View: scroll back so column A is visible again (topLeftCell reset to
# default) and select A15, mirroring the saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A15").Select() | Out-Null
